$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 59
$ws.Range("A59").Value = "ukb51139_subset.csv"
$ws.Range("B59").Value = "5602 x 1081"
$ws.Range("C59").Value = "all"
$ws.Range("D59").Value = "no events"
$ws.Range("E59").Value = "> 160/100"
$ws.Range("F59").Value = "zscore"
$ws.Range("G59").Value = "median"
$ws.Range("H59").Value = "none"
$ws.Range("I59").Value = 25
$ws.Range("L59").Value = "98.9 & 86.8"
$ws.Range("M59").Value = "73.1 & 54.6"
$ws.Range("N59").Value = 20
$ws.Range("O59").Value = 3.85
$ws.Range("P59").Value = "feature selection 0.95"

# Row 60
$ws.Range("A60").Value = "ukb51139_subset.csv"
$ws.Range("B60").Value = "5602 x 1081"
$ws.Range("C60").Value = "all"
$ws.Range("D60").Value = "no events"
$ws.Range("E60").Value = "> 160/100"
$ws.Range("F60").Value = "zscore"
$ws.Range("G60").Value = "median"
$ws.Range("H60").Value = "none"
$ws.Range("I60").Value = 25
$ws.Range("L60").Value = "bad"
$ws.Range("M60").Value = "bad"
$ws.Range("N60").Value = 11
$ws.Range("O60").Value = 91.1
$ws.Range("P60").Value = "feature selection 0.5"
